$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.701.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.095.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5155"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4392"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09270"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.092.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.301"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.753"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06672"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.730.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.324"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.334.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "

$ws.Range("E28").Value = "  -5.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.133"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.41%  "

$ws.Range("E32").Value = "  -2.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.171"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.938"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.202"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02571"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06697"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6868"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2228"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.306"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6671"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.324"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.629"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000353"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3289"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
